# Phase 2 & 3: refresh the exported MRO-tracking data with a new record
# and drop the extra sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "time"/"date" columns swap labels (B <-> C)
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "time"

# Row 2: replace the data with the new record
$ws.Range("A2").Value = "5dd69de6fa189d2ae048845c"
$ws.Range("B2").Value = "Thu Nov 21 2019 00:00:00 GMT+0530 (India Standard Time)"
$ws.Range("C2").Value = "19:45:09"
$ws.Range("D2").Value = '{"_id":"5dab54277715631de4b15d13","dist":"Srikakulam","mandal":"Veeraghattam","mroPhone":"7995995849","hasTelegram":true,"hasWhatsApp":false,"mroName":"#N/a","__v":0}'
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = "19:53"
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = "19:53"
$ws.Range("I2").Value = 0

# Remove the now-stale sample rows 3-5
$ws.Rows("3:5").Delete()
